$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Value = "11_02_2024"
$ws.Range("H2").Value = 3215
$ws.Range("H3").Value = 2385
$ws.Range("H4").Value = 3682
$ws.Range("H5").Value = 6931

$ws.Range("H6").Select()
